$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("nonlinear")

# Fix the "short" (column A) labels: drop the leading zero in the
# MSE01..MSE09 labels, and drop the space in "alpha 1"/"alpha 2".
$ws.Range("A16").Value = "MSE1"
$ws.Range("A17").Value = "MSE2"
$ws.Range("A18").Value = "MSE3"
$ws.Range("A19").Value = "MSE4"
$ws.Range("A20").Value = "MSE5"
$ws.Range("A21").Value = "MSE6"
$ws.Range("A23").Value = "MSE8"
$ws.Range("A22").Value = "MSE7"
$ws.Range("A24").Value = "MSE9"
$ws.Range("A8").Value = "alpha2"
$ws.Range("A7").Value = "alpha1"

# Make the sheet active and move the selection to A7, matching the
# saved cursor position in the workbook.
$ws.Activate()
$ws.Range("A7").Select() | Out-Null
